$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: introduce new shared strings in the exact order the original author typed them ---
# (so the rebuilt sharedStrings.xml table gets the same index assignment as the target file)
$ws.Range("A143").Value = '##########################################################################################################'
$ws.Range("B145").Value = 'Deck: Wilderness'
$ws.Range("B147").Value = 'Artifacts: '
$ws.Range("G151").Value = 'Creatures'
$ws.Range("B151").Value = 'Ash Spawn'
$ws.Range("B182").Value = 'charus'
$ws.Range("D182").Value = 'archer'
$ws.Range("G182").Value = 'creatures'
$ws.Range("B152").Value = 'troll'
$ws.Range("B153").Value = 'frost troll'
$ws.Range("B166").Value = 'ice wraith'
$ws.Range("D166").Value = 'ice mage'
$ws.Range("B155").Value = 'wolf'
$ws.Range("B160").Value = 'ice wolf'
$ws.Range("B156").Value = 'mammoth'
$ws.Range("B157").Value = 'sabrecat'
$ws.Range("B158").Value = 'skeever'
$ws.Range("B167").Value = 'spriggan'
$ws.Range("B168").Value = 'spriggan matron'
$ws.Range("B169").Value = 'spriggan earth mother'
$ws.Range("D167").Value = 'healer'
$ws.Range("B170").Value = 'wispmother'
$ws.Range("B159").Value = 'werewolf'
$ws.Range("B184").Value = 'falmer archer'
$ws.Range("B171").Value = 'hagraven'
$ws.Range("B172").Value = 'glenmorin witch'
$ws.Range("B185").Value = 'charus hunter'
$ws.Range("B154").Value = 'udefrykte(named troll)'
$ws.Range("B161").Value = 'fierce sabrecat'
$ws.Range("B162").Value = 'giant'

# --- Phase 2: fill in remaining cells (numbers + already-existing shared strings / repeats) ---
$ws.Range("A149").Value = 'Warriors'

$ws.Range("A150").Value = 'No.'
$ws.Range("B150").Value = 'Name'
$ws.Range("C150").Value = 'Strength'
$ws.Range("D150").Value = 'Subtype'
$ws.Range("E150").Value = 'Created'
$ws.Range("F150").Value = 'Hero'
$ws.Range("G150").Value = 'Race'
$ws.Range("H150").Value = 'Ability'

$ws.Range("A151").Value = 1
$ws.Range("C151").Value = 1
$ws.Range("D151").Value = 'light'
$ws.Range("E151").Value = 'yes'

$ws.Range("A152").Value = 2
$ws.Range("C152").Value = 3
$ws.Range("D152").Value = 'light'
$ws.Range("E152").Value = 'yes'

$ws.Range("A153").Value = 3
$ws.Range("C153").Value = 4
$ws.Range("D153").Value = 'heavy'
$ws.Range("E153").Value = 'yes'

$ws.Range("A154").Value = 4
$ws.Range("C154").Value = 5
$ws.Range("D154").Value = 'heavy'
$ws.Range("E154").Value = 'yes'

$ws.Range("A155").Value = 5
$ws.Range("C155").Value = 2
$ws.Range("D155").Value = 'light'
$ws.Range("E155").Value = 'yes'

$ws.Range("A156").Value = 6
$ws.Range("C156").Value = 5
$ws.Range("D156").Value = 'heavy'
$ws.Range("E156").Value = 'yes'

$ws.Range("A157").Value = 7
$ws.Range("C157").Value = 6
$ws.Range("D157").Value = 'heavy'
$ws.Range("E157").Value = 'yes'

$ws.Range("A158").Value = 8
$ws.Range("C158").Value = 1
$ws.Range("D158").Value = 'light'
$ws.Range("E158").Value = 'yes'

$ws.Range("A159").Value = 9
$ws.Range("C159").Value = 8
$ws.Range("D159").Value = 'heavy'
$ws.Range("E159").Value = 'yes'

$ws.Range("A160").Value = 10
$ws.Range("C160").Value = 3
$ws.Range("D160").Value = 'light'
$ws.Range("E160").Value = 'yes'

$ws.Range("A161").Value = 11
$ws.Range("C161").Value = 7
$ws.Range("D161").Value = 'heavy'
$ws.Range("E161").Value = 'yes'

$ws.Range("A162").Value = 12
$ws.Range("C162").Value = 9
$ws.Range("D162").Value = 'heavy'
$ws.Range("E162").Value = 'yes'

$ws.Range("A164").Value = 'Mages'

$ws.Range("A165").Value = 'No.'
$ws.Range("B165").Value = 'Name'
$ws.Range("C165").Value = 'Strength'
$ws.Range("D165").Value = 'Subtype'
$ws.Range("E165").Value = 'Created'
$ws.Range("F165").Value = 'Hero'
$ws.Range("G165").Value = 'Race'
$ws.Range("H165").Value = 'Ability'

$ws.Range("A166").Value = 1
$ws.Range("C166").Value = 6

$ws.Range("A167").Value = 2
$ws.Range("C167").Value = 4

$ws.Range("A168").Value = 3
$ws.Range("C168").Value = 7
$ws.Range("D168").Value = 'healer'

$ws.Range("A169").Value = 4
$ws.Range("C169").Value = 10
$ws.Range("D169").Value = 'healer'
$ws.Range("F169").Value = 'yes'

$ws.Range("A170").Value = 5
$ws.Range("C170").Value = 11
$ws.Range("D170").Value = 'lightning mage'
$ws.Range("F170").Value = 'yes'

$ws.Range("A171").Value = 6
$ws.Range("C171").Value = 5
$ws.Range("D171").Value = 'fire mage'

$ws.Range("A172").Value = 7
$ws.Range("C172").Value = 8
$ws.Range("D172").Value = 'fire mage'

$ws.Range("A180").Value = 'Shadow'

$ws.Range("A181").Value = 'No.'
$ws.Range("B181").Value = 'Name'
$ws.Range("C181").Value = 'Strength'
$ws.Range("D181").Value = 'Subtype'
$ws.Range("E181").Value = 'Created'
$ws.Range("F181").Value = 'Hero'
$ws.Range("G181").Value = 'Race'
$ws.Range("H181").Value = 'Ability'

$ws.Range("A182").Value = 1
$ws.Range("C182").Value = 3

$ws.Range("A183").Value = 2

$ws.Range("A184").Value = 3

$ws.Range("A185").Value = 4

$ws.Range("A186").Value = 5

$ws.Range("A187").Value = 6

$ws.Range("A188").Value = 7

$ws.Range("A189").Value = 8

# Update the view selection to match the post-edit state
$ws.Range("G157").Select()
